$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Home, About, Contact"
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Home, About, Contact*") {
        $targetIndex = $i
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    Write-Host "ERROR: could not find paragraph containing 'Home, About, Contact'"
} else {
    $targetPara = $d.Paragraphs.Item($targetIndex)

    # Replace its text with "HOME" (keeps the paragraph's own formatting/pPr)
    $targetPara.Range.Text = "HOME"

    # Insert a brand new paragraph right after it, inheriting the same
    # paragraph formatting (spacing line=240/auto), then give it the text
    # "ABOUT – CONTACT" (en dash U+2013).
    $targetPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "ABOUT " + [char]0x2013 + " CONTACT"

    Write-Host ("Paragraph " + $targetIndex + " now reads: [" + $targetPara.Range.Text + "]")
    Write-Host ("Paragraph " + ($targetIndex + 1) + " now reads: [" + $newPara.Range.Text + "]")
}
